$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '41.758.51'
$ws.Range("E2").Value = '  +0.57%  '

# Row 3
$ws.Range("D3").Value = '2.478.37'
$ws.Range("E3").Value = '  +0.42%  '

# Row 4
$ws.Range("D4").Value = '''1.00'
$ws.Range("E4").Value = '  -0.03%  '

# Row 5
$ws.Range("D5").Value = '''319.45'
$ws.Range("E5").Value = '  +1.50%  '

# Row 6
$ws.Range("D6").Value = '''93.43'
$ws.Range("E6").Value = '  +1.82%  '

# Row 7
$ws.Range("E7").Value = '  +0.78%  '

# Row 8
$ws.Range("E8").Value = '  +0.00%  '

# Row 9
$ws.Range("E9").Value = '  +1.49%  '

# Row 10
$ws.Range("E10").Value = '  +10.15%  '

# Row 11
$ws.Range("D11").Value = '''33.31'
$ws.Range("E11").Value = '  +3.40%  '

# Row 12
$ws.Range("E12").Value = '  +0.69%  '

# Row 13
$ws.Range("D13").Value = '2.860.72'
$ws.Range("E13").Value = '  +0.42%  '

# Row 14
$ws.Range("D14").Value = '''6.91'
$ws.Range("E14").Value = '  +1.20%  '

# Row 15
$ws.Range("D15").Value = '''15.80'
$ws.Range("E15").Value = '  -0.66%  '

# Row 16
$ws.Range("D16").Value = '2.471.88'
$ws.Range("E16").Value = '  +0.04%  '

# Row 17
$ws.Range("D17").Value = '''0.794'
$ws.Range("E17").Value = '  +2.91%  '

# Row 18
$ws.Range("D18").Value = '41.749.15'
$ws.Range("E18").Value = '  +0.56%  '

# Row 19
$ws.Range("E19").Value = '  +0.36%  '

# Row 20
$ws.Range("E20").Value = '  +1.26%  '

# Row 21
$ws.Range("D21").Value = '''71.30'
$ws.Range("E21").Value = '  +0.14%  '

# Row 22
$ws.Range("E22").Value = '  +2.70%  '

# Row 23
$ws.Range("D23").Value = '''239.97'
$ws.Range("E23").Value = '  +1.85%  '

# Row 24
$ws.Range("D24").Value = '''2.75'
$ws.Range("E24").Value = '  +1.24%  '

# Row 25
$ws.Range("E25").Value = '  +2.51%  '

# Row 26
$ws.Range("E26").Value = '  +0.06%  '

# Row 27
$ws.Range("D27").Value = '''24.76'
$ws.Range("E27").Value = '  +0.56%  '

# Row 28
$ws.Range("E28").Value = '  +1.37%  '

# Row 29
$ws.Range("D29").Value = '''9.84'
$ws.Range("E29").Value = '  +1.74%  '

# Row 30
$ws.Range("D30").Value = '''36.23'

# Row 31
$ws.Range("D31").Value = '''158.24'
$ws.Range("E31").Value = '  +1.68%  '

# Row 32
$ws.Range("D32").Value = '''5.53'
$ws.Range("E32").Value = '  +1.91%  '

# Row 33
$ws.Range("E33").Value = '  -0.19%  '

# Row 34
$ws.Range("E34").Value = '  +0.56%  '

# Row 35
$ws.Range("E35").Value = '  +1.56%  '

# Row 36
$ws.Range("D36").Value = '''17.47'
$ws.Range("E36").Value = '  +0.84%  '

# Row 37
$ws.Range("E37").Value = '  +6.68%  '

# Row 38
$ws.Range("E38").Value = '  +2.55%  '

# Row 39
$ws.Range("E39").Value = '  +1.81%  '

# Row 40
$ws.Range("E40").Value = '  +0.45%  '

# Row 41
$ws.Range("D41").Value = '''4.06'
$ws.Range("E41").Value = '  +1.01%  '

# Row 42
$ws.Range("E42").Value = '  +11.24%  '

# Row 43
$ws.Range("D43").Value = '1.993.41'
$ws.Range("E43").Value = '  +2.55%  '

# Row 44
$ws.Range("B44").Value = 'VeChain'
$ws.Range("C44").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D44").Value = '''0.0286'
$ws.Range("E44").Value = '  +1.11%  '

# Row 45
$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").Value = '''19.06'
$ws.Range("E45").Value = '  +3.44%  '

# Row 46
$ws.Range("E46").Value = '  +2.84%  '

# Row 47
$ws.Range("E47").Value = '  +5.27%  '

# Row 48
$ws.Range("D48").Value = '2.717.95'
$ws.Range("E48").Value = '  +0.41%  '

# Row 49
$ws.Range("D49").Value = '''97.65'
$ws.Range("E49").Value = '  +0.93%  '

# Row 50
$ws.Range("D50").Value = '''74.40'
$ws.Range("E50").Value = '  +3.59%  '

# Row 51
$ws.Range("D51").Value = '''67.37'
$ws.Range("E51").Value = '  +0.59%  '
